# QA update: fix leading-space typo in the "Primary Network Interface" label.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "''Primary Network Interface'"
